$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A70").Value = 46030
$ws.Range("B70").Value = 707
$ws.Range("C70").Value = 692
$ws.Range("D70").Value = 15

$ws.Range("A70:D70").Select()
